$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1026.9412
$ws.Cells.Item(15, 9).Value = 1026.9412
$ws.Cells.Item(15, 11).Value = 3080.8236
$ws.Cells.Item(15, 13).Value = -2911.8236
$ws.Cells.Item(33, 8).Value = 1002.25
$ws.Cells.Item(33, 9).Value = 95.07692
$ws.Cells.Item(33, 10).Value = 2687
$ws.Cells.Item(33, 11).Value = 95.07692
$ws.Cells.Item(33, 12).Value = 2687
$ws.Cells.Item(33, 13).Value = 133.92308
$ws.Cells.Item(33, 14).Value = -3145
$ws.Cells.Item(98, 8).Value = 1070.8077
$ws.Cells.Item(98, 9).Value = 1015.7059
$ws.Cells.Item(98, 10).Value = 1174.8889
$ws.Cells.Item(98, 11).Value = 1015.7059
$ws.Cells.Item(98, 12).Value = 1174.8889
$ws.Cells.Item(98, 13).Value = 482.2941
$ws.Cells.Item(98, 14).Value = -4170.8889
$ws.Cells.Item(112, 8).Value = 782.3077
$ws.Cells.Item(112, 10).Value = 995.7143
$ws.Cells.Item(112, 12).Value = 2987.1429
$ws.Cells.Item(112, 14).Value = -5203.1429
$ws.Cells.Item(122, 8).Value = 1070.8077
$ws.Cells.Item(122, 9).Value = 1015.7059
$ws.Cells.Item(122, 10).Value = 1174.8889
$ws.Cells.Item(122, 11).Value = 3047.1177
$ws.Cells.Item(122, 12).Value = 3524.6667
$ws.Cells.Item(122, 13).Value = -597.1177000000002
$ws.Cells.Item(122, 14).Value = -8424.6667
$ws.Cells.Item(132, 8).Value = 14714459
$ws.Cells.Item(132, 9).Value = 15633988
$ws.Cells.Item(132, 10).Value = 1999
$ws.Cells.Item(132, 11).Value = 46901964
$ws.Cells.Item(132, 12).Value = 5997
$ws.Cells.Item(132, 13).Value = -46899434
$ws.Cells.Item(132, 14).Value = -11057
$ws.Cells.Item(138, 8).Value = 5869.343
$ws.Cells.Item(138, 9).Value = 1735.8948
$ws.Cells.Item(138, 10).Value = 10777.8125
$ws.Cells.Item(138, 11).Value = 5207.6844
$ws.Cells.Item(138, 12).Value = 32333.4375
$ws.Cells.Item(138, 13).Value = -67.6844000000001
$ws.Cells.Item(138, 14).Value = -42613.4375
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 2358.5
$ws.Cells.Item(61, 9).Value = 1631.1
$ws.Cells.Item(61, 10).Value = 2638.2693
$ws.Cells.Item(61, 11).Value = 1631.1
$ws.Cells.Item(61, 12).Value = 2638.2693
$ws.Cells.Item(61, 13).Value = -1419.1
$ws.Cells.Item(61, 14).Value = -3062.2693
$ws.Cells.Item(69, 8).Value = 31361.777
$ws.Cells.Item(69, 10).Value = 31361.777
$ws.Cells.Item(69, 12).Value = 31361.777
$ws.Cells.Item(69, 14).Value = -32859.777
$ws.Cells.Item(72, 8).Value = 31361.777
$ws.Cells.Item(72, 10).Value = 31361.777
$ws.Cells.Item(72, 12).Value = 94085.33099999999
$ws.Cells.Item(72, 14).Value = -101573.331
$ws.Cells.Item(74, 8).Value = 1829.7028
$ws.Cells.Item(74, 9).Value = 2111.0833
$ws.Cells.Item(74, 10).Value = 1694.64
$ws.Cells.Item(74, 11).Value = 2111.0833
$ws.Cells.Item(74, 12).Value = 1694.64
$ws.Cells.Item(74, 13).Value = -1237.0833
$ws.Cells.Item(74, 14).Value = -3442.64
$ws.Cells.Item(77, 8).Value = 1829.7028
$ws.Cells.Item(77, 9).Value = 2111.0833
$ws.Cells.Item(77, 10).Value = 1694.64
$ws.Cells.Item(77, 11).Value = 10555.4165
$ws.Cells.Item(77, 12).Value = 8473.200000000001
$ws.Cells.Item(77, 13).Value = -6187.416499999999
$ws.Cells.Item(77, 14).Value = -17209.2
$ws.Cells.Item(122, 8).Value = 2060.9644
$ws.Cells.Item(122, 9).Value = 1856.0454
$ws.Cells.Item(122, 11).Value = 5568.1362
$ws.Cells.Item(122, 13).Value = -3118.1362
$ws.Cells.Item(132, 8).Value = 6030.75
$ws.Cells.Item(132, 9).Value = 6541.1665
$ws.Cells.Item(132, 10).Value = 4499.5
$ws.Cells.Item(132, 11).Value = 19623.4995
$ws.Cells.Item(132, 12).Value = 13498.5
$ws.Cells.Item(132, 13).Value = -17093.4995
$ws.Cells.Item(132, 14).Value = -18558.5
$ws.Cells.Item(136, 8).Value = 2358.5
$ws.Cells.Item(136, 9).Value = 1631.1
$ws.Cells.Item(136, 10).Value = 2638.2693
$ws.Cells.Item(136, 11).Value = 4893.299999999999
$ws.Cells.Item(136, 12).Value = 7914.8079
$ws.Cells.Item(136, 13).Value = -2343.299999999999
$ws.Cells.Item(136, 14).Value = -13014.8079
$ws.Cells.Item(139, 8).Value = 61968
$ws.Cells.Item(139, 10).Value = 61968
$ws.Cells.Item(139, 12).Value = 61968
$ws.Cells.Item(139, 14).Value = -72248
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(75, 8).Value = 7699
$ws.Cells.Item(75, 9).Value = 8438.799999999999
$ws.Cells.Item(75, 10).Value = 4000
$ws.Cells.Item(75, 11).Value = 8438.799999999999
$ws.Cells.Item(75, 12).Value = 4000
$ws.Cells.Item(75, 13).Value = -7502.799999999999
$ws.Cells.Item(75, 14).Value = -5872
$ws.Cells.Item(78, 8).Value = 7699
$ws.Cells.Item(78, 9).Value = 8438.799999999999
$ws.Cells.Item(78, 10).Value = 4000
$ws.Cells.Item(78, 11).Value = 25316.4
$ws.Cells.Item(78, 12).Value = 12000
$ws.Cells.Item(78, 13).Value = -20636.4
$ws.Cells.Item(78, 14).Value = -21360
$ws.Cells.Item(105, 8).Value = 202030.1
$ws.Cells.Item(105, 9).Value = 126661.25
$ws.Cells.Item(105, 10).Value = 503505.5
$ws.Cells.Item(105, 11).Value = 126661.25
$ws.Cells.Item(105, 12).Value = 503505.5
$ws.Cells.Item(105, 13).Value = -124914.25
$ws.Cells.Item(105, 14).Value = -506999.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 5369
$ws.Cells.Item(132, 9).Value = 7187.2
$ws.Cells.Item(132, 10).Value = 3853.8333
$ws.Cells.Item(132, 11).Value = 21561.6
$ws.Cells.Item(132, 12).Value = 11561.4999
$ws.Cells.Item(132, 13).Value = -19031.6
$ws.Cells.Item(132, 14).Value = -16621.4999
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(18, 8).Value = 303.875
$ws.Cells.Item(18, 9).Value = 254.53847
$ws.Cells.Item(18, 11).Value = 763.61541
$ws.Cells.Item(18, 13).Value = -594.61541
$ws.Cells.Item(37, 8).Value = 1133115.6
$ws.Cells.Item(37, 10).Value = 1133115.6
$ws.Cells.Item(37, 12).Value = 3399346.8
$ws.Cells.Item(37, 14).Value = -3399570.8
$ws.Cells.Item(45, 8).Value = 922.5
$ws.Cells.Item(45, 9).Value = 845
$ws.Cells.Item(45, 10).Value = 1000
$ws.Cells.Item(45, 11).Value = 2535
$ws.Cells.Item(45, 12).Value = 3000
$ws.Cells.Item(45, 13).Value = -2003
$ws.Cells.Item(45, 14).Value = -4064
$ws.Cells.Item(52, 8).Value = 1839.7778
$ws.Cells.Item(52, 10).Value = 1839.7778
$ws.Cells.Item(52, 12).Value = 5519.3334
$ws.Cells.Item(52, 14).Value = -6051.3334
$ws.Cells.Item(54, 8).Value = 0
$ws.Cells.Item(54, 10).Value = 0
$ws.Cells.Item(54, 12).Value = 0
$ws.Cells.Item(54, 14).Value = ""
$ws.Cells.Item(120, 8).Value = 8633.333000000001
$ws.Cells.Item(120, 10).Value = 7950
$ws.Cells.Item(120, 12).Value = 23850
$ws.Cells.Item(120, 14).Value = -33526
$ws.Cells.Item(121, 8).Value = 12658.583
$ws.Cells.Item(121, 9).Value = 9354.75
$ws.Cells.Item(121, 10).Value = 14310.5
$ws.Cells.Item(121, 11).Value = 28064.25
$ws.Cells.Item(121, 12).Value = 42931.5
$ws.Cells.Item(121, 13).Value = -26754.25
$ws.Cells.Item(121, 14).Value = -45551.5
$ws.Cells.Item(123, 8).Value = 2860.8
$ws.Cells.Item(123, 9).Value = 1534.8334
$ws.Cells.Item(123, 10).Value = 4849.75
$ws.Cells.Item(123, 11).Value = 4604.5002
$ws.Cells.Item(123, 12).Value = 14549.25
$ws.Cells.Item(123, 13).Value = -2154.5002
$ws.Cells.Item(123, 14).Value = -19449.25
$ws.Cells.Item(124, 8).Value = 4552
$ws.Cells.Item(124, 9).Value = 3000
$ws.Cells.Item(124, 10).Value = 4940
$ws.Cells.Item(124, 11).Value = 9000
$ws.Cells.Item(124, 12).Value = 14820
$ws.Cells.Item(124, 13).Value = -4090
$ws.Cells.Item(124, 14).Value = -24640
$ws.Cells.Item(129, 8).Value = 7813959
$ws.Cells.Item(129, 9).Value = 19231244
$ws.Cells.Item(129, 10).Value = 2133.3157
$ws.Cells.Item(129, 11).Value = 57693732
$ws.Cells.Item(129, 12).Value = 6399.9471
$ws.Cells.Item(129, 13).Value = -57688732
$ws.Cells.Item(129, 14).Value = -16399.9471
$ws.Cells.Item(130, 8).Value = 2080
$ws.Cells.Item(130, 9).Value = 0
$ws.Cells.Item(130, 10).Value = 2080
$ws.Cells.Item(130, 11).Value = 0
$ws.Cells.Item(130, 12).Value = 6240
$ws.Cells.Item(130, 13).Value = ""
$ws.Cells.Item(130, 14).Value = -16280
$ws.Cells.Item(131, 8).Value = 1456.4584
$ws.Cells.Item(131, 9).Value = 318.75
$ws.Cells.Item(131, 10).Value = 1559.8864
$ws.Cells.Item(131, 11).Value = 956.25
$ws.Cells.Item(131, 12).Value = 4679.6592
$ws.Cells.Item(131, 13).Value = 4083.75
$ws.Cells.Item(131, 14).Value = -14759.6592
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 4130.2856
$ws.Cells.Item(132, 9).Value = 3801.3333
$ws.Cells.Item(132, 10).Value = 4722.4
$ws.Cells.Item(132, 11).Value = 11403.9999
$ws.Cells.Item(132, 12).Value = 14167.2
$ws.Cells.Item(132, 13).Value = -8873.999899999999
$ws.Cells.Item(132, 14).Value = -19227.2
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 689.1905
$ws.Cells.Item(55, 9).Value = 392.35715
$ws.Cells.Item(55, 10).Value = 1282.8572
$ws.Cells.Item(55, 11).Value = 392.35715
$ws.Cells.Item(55, 12).Value = 1282.8572
$ws.Cells.Item(55, 13).Value = -219.35715
$ws.Cells.Item(55, 14).Value = -1628.8572
$ws.Cells.Item(123, 8).Value = 24964.312
$ws.Cells.Item(123, 10).Value = 24964.312
$ws.Cells.Item(123, 12).Value = 24964.312
$ws.Cells.Item(123, 14).Value = -34764.31200000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 766.54285
$ws.Cells.Item(136, 9).Value = 545.3182
$ws.Cells.Item(136, 11).Value = 1635.9546
$ws.Cells.Item(136, 13).Value = 914.0454
